$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the empty column G (shifts H,I -> G,H and M -> L)
$ws.Columns.Item(7).Delete()

# Swap columns A and B (data + styles): A currently has Dates, B has Model #.
# Target: A should have Model #, B should have Dates.
$ws.Columns.Item(2).Cut()
$ws.Columns.Item(1).Insert()

# Fix the selection to match the target view
$ws.Range("J17").Select()
